# Edit workbook per commit: "changed shfitish to ngrams, revamped arcane ngrams,
# added dash toggle to sniping, arcane ellipsis"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 header "key" -> five apostrophes ('''''), stored as text with a
# leading quote-prefix (Excel eats one leading apostrophe as the
# quote-prefix marker, so we type six to end up with five stored chars).
$ws.Range("A1").Value = "''''''"

# New cell H8 gets a "~" (dash/tilde toggle entry added to the sniping block)
$ws.Range("H8").Value = "~"

# L13 changes from "~" to "-"
$ws.Range("L13").Value = "-"

# Move the active selection to K13 (matches the saved cursor position)
$ws.Range("K13").Select() | Out-Null
